$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Cell row=1 col=1: "58 x 75" -> "59 x 23"
$cell = $t.Cell(1, 1)
$cell.Range.Text = "59 x 23`v  2    3`v  ----`v5|    |`v9|    |"

# Cell row=1 col=2: "41 x 79" -> "24 x 99"
$cell = $t.Cell(1, 2)
$cell.Range.Text = "24 x 99`v  9    9`v  ----`v2|    |`v4|    |"

# Cell row=1 col=3: "95 x 28" -> "74 x 40"
$cell = $t.Cell(1, 3)
$cell.Range.Text = "74 x 40`v  4    0`v  ----`v7|    |`v4|    |"

# Cell row=2 col=1: "55 x 27" -> "16 x 18"
$cell = $t.Cell(2, 1)
$cell.Range.Text = "16 x 18`v  1    8`v  ----`v1|    |`v6|    |"

# Cell row=2 col=2: "64 x 70" -> "77 x 58"
$cell = $t.Cell(2, 2)
$cell.Range.Text = "77 x 58`v  5    8`v  ----`v7|    |`v7|    |"

# Cell row=2 col=3: "94 x 19" -> "80 x 77"
$cell = $t.Cell(2, 3)
$cell.Range.Text = "80 x 77`v  7    7`v  ----`v8|    |`v0|    |"

# Cell row=3 col=1: "61 x 17" -> "32 x 17"
$cell = $t.Cell(3, 1)
$cell.Range.Text = "32 x 17`v  1    7`v  ----`v3|    |`v2|    |"

# Cell row=3 col=2: "45 x 47" -> "65 x 19"
$cell = $t.Cell(3, 2)
$cell.Range.Text = "65 x 19`v  1    9`v  ----`v6|    |`v5|    |"

# Cell row=3 col=3: "34 x 52" -> "63 x 50"
$cell = $t.Cell(3, 3)
$cell.Range.Text = "63 x 50`v  5    0`v  ----`v6|    |`v3|    |"

# Cell row=4 col=1: "83 x 85" -> "79 x 54"
$cell = $t.Cell(4, 1)
$cell.Range.Text = "79 x 54`v  5    4`v  ----`v7|    |`v9|    |"

# Cell row=4 col=2: "91 x 30" -> "79 x 39"
$cell = $t.Cell(4, 2)
$cell.Range.Text = "79 x 39`v  3    9`v  ----`v7|    |`v9|    |"

# Cell row=4 col=3: "27 x 82" -> "63 x 81"
$cell = $t.Cell(4, 3)
$cell.Range.Text = "63 x 81`v  8    1`v  ----`v6|    |`v3|    |"

# Cell row=5 col=1: "67 x 83" -> "86 x 66"
$cell = $t.Cell(5, 1)
$cell.Range.Text = "86 x 66`v  6    6`v  ----`v8|    |`v6|    |"

# Cell row=5 col=2: "74 x 62" -> "45 x 15"
$cell = $t.Cell(5, 2)
$cell.Range.Text = "45 x 15`v  1    5`v  ----`v4|    |`v5|    |"

# Cell row=5 col=3: "17 x 99" -> "87 x 12"
$cell = $t.Cell(5, 3)
$cell.Range.Text = "87 x 12`v  1    2`v  ----`v8|    |`v7|    |"
